$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing quarterly job counts (values revised with newer data)
$ws.Range("B2").Value = 418
$ws.Range("B4").Value = 117
$ws.Range("B5").Value = 140
$ws.Range("B6").Value = 778
$ws.Range("B7").Value = 138
$ws.Range("B8").Value = 144
$ws.Range("B9").Value = 193
$ws.Range("B10").Value = 1103
$ws.Range("B11").Value = 225
$ws.Range("B12").Value = 195
$ws.Range("B13").Value = 291
$ws.Range("B14").Value = 1982
$ws.Range("B15").Value = 364
$ws.Range("B16").Value = 392
$ws.Range("B17").Value = 369
$ws.Range("B19").Value = 573
$ws.Range("B20").Value = 488
$ws.Range("B21").Value = 621
$ws.Range("B22").Value = 1685
$ws.Range("B23").Value = 1407
$ws.Range("B24").Value = 639

# Add new 2020 Q4 data row
$ws.Range("A25").Value = "2020-Q4"
$ws.Range("B25").Value = 542

# Match number formatting used by the rest of the "Jobs Filled" column
$ws.Range("B24").Copy()
$ws.Range("B25").PasteSpecial(-4122)
